$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the data block (row 311), pushing the
# existing rows 311:361 down to 314:364. Excel's Insert() copies the
# formatting of the row above, which preserves the date-style (s="2")
# on column D for the shifted rows as well as the newly inserted ones.
$ws.Rows("311:313").Insert()

# Populate the 3 newly inserted rows with the new weekly report
# (fecha 2021-11-04 / serial 44504), following the same
# Maduro / Pintón / Primera Pintón pattern used throughout the sheet.

# Row 311 - Maduro
$ws.Range("A311").Value = 11
$ws.Range("B311").Value = "Vega Monumental Concepción"
$ws.Range("C311").Value = "Bíobío"
$ws.Range("D311").Value = 44504
$ws.Range("E311").Value = 8
$ws.Range("F311").Value = "Fruta"
$ws.Range("G311").Value = 100108
$ws.Range("H311").Value = "Tropicales y subtropicales"
$ws.Range("I311").Value = 100108006
$ws.Range("J311").Value = "Plátano"
$ws.Range("K311").Value = "Sin especificar"
$ws.Range("L311").Value = "Maduro"
$ws.Range("M311").Value = 100
$ws.Range("N311").Value = 13000
$ws.Range("O311").Value = 13000
$ws.Range("P311").Value = 13000
$ws.Range("Q311").Value = "$/caja 20 kilos"
$ws.Range("R311").Value = "Ecuador"
$ws.Range("S311").Value = 650
$ws.Range("T311").Value = 20

# Row 312 - Pintón
$ws.Range("A312").Value = 11
$ws.Range("B312").Value = "Vega Monumental Concepción"
$ws.Range("C312").Value = "Bíobío"
$ws.Range("D312").Value = 44504
$ws.Range("E312").Value = 8
$ws.Range("F312").Value = "Fruta"
$ws.Range("G312").Value = 100108
$ws.Range("H312").Value = "Tropicales y subtropicales"
$ws.Range("I312").Value = 100108006
$ws.Range("J312").Value = "Plátano"
$ws.Range("K312").Value = "Sin especificar"
$ws.Range("L312").Value = "Pintón"
$ws.Range("M312").Value = 300
$ws.Range("N312").Value = 14000
$ws.Range("O312").Value = 14000
$ws.Range("P312").Value = 14000
$ws.Range("Q312").Value = "$/caja 20 kilos"
$ws.Range("R312").Value = "Ecuador"
$ws.Range("S312").Value = 700
$ws.Range("T312").Value = 20

# Row 313 - Primera Pintón
$ws.Range("A313").Value = 11
$ws.Range("B313").Value = "Vega Monumental Concepción"
$ws.Range("C313").Value = "Bíobío"
$ws.Range("D313").Value = 44504
$ws.Range("E313").Value = 8
$ws.Range("F313").Value = "Fruta"
$ws.Range("G313").Value = 100108
$ws.Range("H313").Value = "Tropicales y subtropicales"
$ws.Range("I313").Value = 100108006
$ws.Range("J313").Value = "Plátano"
$ws.Range("K313").Value = "Sin especificar"
$ws.Range("L313").Value = "Primera Pintón"
$ws.Range("M313").Value = 300
$ws.Range("N313").Value = 16000
$ws.Range("O313").Value = 16000
$ws.Range("P313").Value = 16000
$ws.Range("Q313").Value = "$/caja 20 kilos"
$ws.Range("R313").Value = "Ecuador"
$ws.Range("S313").Value = 800
$ws.Range("T313").Value = 20
